$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the username/password pairs down by one row (B2:C9 -> B3:C10),
# leaving the id column (A) untouched, and add a new id for the row
# that is pushed out to row 10.
for ($r = 9; $r -ge 2; $r--) {
    $bval = $ws.Cells.Item($r, 2).Value()
    $cval = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($r + 1, 2).Value = $bval
    $ws.Cells.Item($r + 1, 3).Value = $cval
}

# New id for the row that moved from 9 to 10.
$ws.Cells.Item(10, 1).Value = 10

# Fill in the freed-up row 2 with the new admin account.
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "admin"
$ws.Cells.Item(2, 3).Value = "admin"

# Match the final selection left behind in the sheet.
$ws.Range("C6").Select()
